$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new diary entries: rows 21-28 --------------------------------
# Columns: A=Data, B=Entrada turno1, C=Saida turno1, D=Entrada turno2,
#          E=Saida turno2, F=Horas no dia, G=Total acumulado, H=Atividade

# Copy the formatting of the last existing data row (20) down into the
# new rows so the new cells pick up the same number formats / styles.
$ws.Range("A20:H20").Copy()
$ws.Range("A21:H27").PasteSpecial(-4122)
# Row 28 is only partially filled in (date + first "entrada" time), so
# only copy formatting for those two columns.
$ws.Range("A20:B20").Copy()
$ws.Range("A28:B28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$rows = @(
    @{ Row=21; A=43892; B=0.375;               C=0.45833333333333331; D=0.79166666666666663; E=0.95833333333333337; F=6; H='Finalização do cartão de vacina e começo do PDF da ficha de atendimento' },
    @{ Row=22; A=43893; B=0.375;               C=0.45833333333333331; D=0.625;                E=0.70833333333333337; F=4; H='Finalização do PDF ficha de atendimento, mudança de data para formato brasileiro' },
    @{ Row=23; A=43895; B=0.375;               C=0.5;                 D=0.54166666666666663; E=0.625;                F=5; H='Reunião com o professor Leonardo e Personalização da interface gráfica 2.0' },
    @{ Row=24; A=43896; B=0.375;               C=0.45833333333333331; D=0.54166666666666663; E=0.70833333333333337; F=6; H='Personalização da interface gráfica (Quiron 2.0)' },
    @{ Row=25; A=43899; B=0.375;               C=0.45833333333333331; D=0.79166666666666663; E=0.875;                F=4; H='Personalização da interface gráfica (Quiron 2.0)' },
    @{ Row=26; A=43900; B=0.33333333333333331; C=0.45833333333333331; D=0.79166666666666663; E=0.91666666666666663; F=6; H='Personalização da interface gráfica (Quiron 2.0) e cabeçalho em cada tela' },
    @{ Row=27; A=43901; B=0.375;               C=0.45833333333333331; D=0.625;                E=0.70833333333333337; F=4; H='Tentativa de criptografia e exportação do banco de dados na enfermaria' }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $prev = $row - 1
    $ws.Cells.Item($row, 7).Formula = "=G$prev+F$row"
    $ws.Cells.Item($row, 8).Value = $r.H
}

# Row 28: still being filled in, only date + entrada turno 1 present so far.
$ws.Cells.Item(28, 1).Value = 43902
$ws.Cells.Item(28, 2).Value = 0.39583333333333331

# --- Update the view to match where the author left off ----------------
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("E32").Select() | Out-Null
